# Second edit to word document
#
# Splits the single paragraph "This is my Word Document" into two
# paragraphs: the original text, and a new second paragraph reading
# "2nd edit to word document" (with "nd" superscripted). The trailing
# _GoBack bookmark is moved from the end of the first paragraph to the
# end of the newly added second paragraph.

$d = $word.ActiveDocument

# The existing _GoBack bookmark currently sits right after
# "This is my Word Document". Remove it now; it will be re-created at
# the end of the new second paragraph once that text exists.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# "This is my Word Document" occupies characters 0-23, with the
# paragraph mark at position 24. Insert a paragraph break followed by
# "2" right before that paragraph mark.
$splitPoint = 24
$rng = $d.Range($splitPoint, $splitPoint)
$rng.InsertAfter("`r2")

# Insert "nd" immediately after the "2" we just added (positions 24+1=25
# is "2", so the new insertion point is 26).
$ndInsertPos = $splitPoint + 2
$ndRng = $d.Range($ndInsertPos, $ndInsertPos)
$ndRng.InsertAfter("nd")

# Format the "nd" we just inserted as superscript.
$ndRng2 = $d.Range($ndInsertPos, $ndInsertPos + 2)
$ndRng2.Font.Superscript = $true

# Insert the remainder of the sentence after "nd".
$restInsertPos = $ndInsertPos + 2
$restRng = $d.Range($restInsertPos, $restInsertPos)
$restRng.InsertAfter(" edit to word document")

# Re-add the _GoBack bookmark at the end of the new second paragraph
# (i.e. right before the document's final paragraph mark).
#
# Note: adding an empty (zero-length) bookmark exactly at
# Content.End-1 (the gap immediately before the very last paragraph
# mark in the document) is mishandled and ends up placed at the start
# of the document instead. Work around this by temporarily inserting a
# placeholder character after the target position (which makes the
# bookmark's position no longer the last gap in the document), adding
# the bookmark there, and then removing the placeholder again.
$beforePlaceholderEnd = $d.Content.End
$placeholderRng = $d.Range($beforePlaceholderEnd, $beforePlaceholderEnd)
$placeholderRng.InsertAfter("X")

# The placeholder character actually lands one position before the
# previous Content.End, since InsertAfter on an empty range at the very
# end of the document inserts just before the trailing paragraph mark.
$bookmarkPos = $beforePlaceholderEnd - 1

$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))

$placeholderCharRng = $d.Range($bookmarkPos, $bookmarkPos + 1)
$placeholderCharRng.Delete()
